$d = $word.ActiveDocument

# 1. Remove the old hidden "_GoBack" bookmark (it will be re-added later at the
#    new location, right before the closing paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. "community based development" -> "community-based development"
$d.Content.Find.Execute("community based development", $true, $false, $false, $false, $false, $true, 1, $false, "community-based development", 2) | Out-Null

# 3. "Ihavestronginterest in" -> "I have strong interest in"
$d.Content.Find.Execute("Ihavestronginterest in", $true, $false, $false, $false, $false, $true, 1, $false, "I have strong interest in", 2) | Out-Null

# 4. "potable safe" -> "portable safe"
$d.Content.Find.Execute("potable safe", $true, $false, $false, $false, $false, $true, 1, $false, "portable safe", 2) | Out-Null

# 5. Merge "...applying the" paragraph with the following "technical and team..." paragraph.
$d.Content.Find.Execute("applying the^p", $true, $false, $false, $false, $false, $true, 1, $false, "applying the ", 2) | Out-Null

# 6. Merge "...my skills and interest that is" paragraph with the following "committed to..." paragraph.
$d.Content.Find.Execute("my skills and interest that is^p", $true, $false, $false, $false, $false, $true, 1, $false, "my skills and interest that is ", 2) | Out-Null

# 7. Add the "_Hlk8652934" bookmark around "Aim excited ... continuous development and growth"
#    (from the "Aim" right after "and hygiene; and " through "growth", excluding the final period).
$r1 = $d.Content
$r1.Find.Execute("Aim excited by the prospect of", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hlkStart = $r1.Start

$r2 = $d.Content
$r2.Find.Execute("continuous development and growth", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hlkEnd = $r2.End

$hlkRange = $d.Range($hlkStart, $hlkEnd)
$d.Bookmarks.Add("_Hlk8652934", $hlkRange) | Out-Null

# 8. Add the "_GoBack" bookmark at the start of the final paragraph through just before the
#    trailing comma of "Sincerely yours,".
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$goBackStart = $lastPara.Range.Start

$r3 = $d.Content
$r3.Find.Execute("Sincerely yours", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackEnd = $r3.End

$goBackRange = $d.Range($goBackStart, $goBackEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
